# Insert a new data row at row 20 (pushes the existing rows 20-87 down to
# 21-88) and populate it with the new weekly record, matching the diff's
# expanded dimension A1:R88.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("20:20").Insert()

$ws.Cells.Item(20, 1).Value  = 11
$ws.Cells.Item(20, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(20, 3).Value  = "Bíobío"
$ws.Cells.Item(20, 4).Value  = 44565
$ws.Cells.Item(20, 5).Value  = 8
$ws.Cells.Item(20, 6).Value  = 100112032
$ws.Cells.Item(20, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(20, 8).Value  = "Sin especificar"
$ws.Cells.Item(20, 9).Value  = "Primera"
$ws.Cells.Item(20, 10).Value = 100
$ws.Cells.Item(20, 11).Value = 7000
$ws.Cells.Item(20, 12).Value = 8000
$ws.Cells.Item(20, 13).Value = 7500
$ws.Cells.Item(20, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(20, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(20, 16).Value = 150
$ws.Cells.Item(20, 17).Value = 50
$ws.Cells.Item(20, 18).Value = "Hortaliza"
